# Fruta / hortaliza, semanal
# Insert a new weekly record at row 12 (pushing the existing rows 12-27
# down to 13-28) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 12; the formatting/contents of the
# old rows 12-27 shift down to 13-28 automatically.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new weekly record.
$ws.Range("A12").Value = 7
$ws.Range("B12").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C12").Value = "Ñuble"
$ws.Range("D12").Value = 44519
$ws.Range("E12").Value = 16
$ws.Range("F12").Value = 100112026
$ws.Range("G12").Value = "Haba"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 6000
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 6500
$ws.Range("N12").Value = "$/saco 25 kilos"
$ws.Range("O12").Value = "Provincia de Diguillín"
$ws.Range("P12").Value = 260
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = "Hortaliza"
